{"js": "const body = context.document.body;\nconst pairs = [\n  [\"2024-04-10 Wednesday\", \"2024-04-11 Thursday\"],\n  [\"482\u00f79=53, 5\", \"537\u00f73=179, 0\"],\n  [\"738\u00f78=92, 2\", \"721\u00f77=103, 0\"],\n  [\"741\u00f72=370, 1\", \"827\u00f75=165, 2\"],\n  [\"931\u00f75=186, 1\", \"789\u00f72=394, 1\"],\n  [\"893\u00f73=297, 2\", \"192\u00f73=64, 0\"],\n  [\"877\u00f79=97, 4\", \"164\u00f77=23, 3\"],\n  [\"501\u00f72=250, 1\", \"986\u00f78=123, 2\"],\n  [\"317\u00f72=158, 1\", \"236\u00f79=26, 2\"],\n  [\"351\u00f79=39, 0\", \"246\u00f76=41, 0\"],\n  [\"419\u00f78=52, 3\", \"839\u00f77=119, 6\"],\n  [\"882\u00f77=126, 0\", \"517\u00f74=129, 1\"],\n  [\"846\u00f79=94, 0\", \"109\u00f75=21, 4\"],\n  [\"338\u00f72=169, 0\", \"988\u00f79=109, 7\"],\n  [\"895\u00f74=223, 3\", \"627\u00f76=104, 3\"],\n  [\"781\u00f78=97, 5\", \"896\u00f79=99, 5\"],\n  [\"257\u00f75=51, 2\", \"233\u00f74=58, 1\"],\n  [\"282\u00f74=70, 2\", \"110\u00f78=13, 6\"],\n  [\"645\u00f73=215, 0\", \"182\u00f79=20, 2\"],\n  [\"496\u00f78=62, 0\", \"140\u00f76=23, 2\"],\n  [\"785\u00f79=87, 2\", \"478\u00f75=95, 3\"],\n  [\"636\u00f77=90, 6\", \"935\u00f77=133, 4\"],\n  [\"771\u00f72=385, 1\", \"693\u00f73=231, 0\"],\n  [\"877\u00f72=438, 1\", \"785\u00f79=87, 2\"],\n  [\"912\u00f79=101, 3\", \"756\u00f78=94, 4\"],\n  [\"914\u00f72=457, 0\", \"680\u00f74=170, 0\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('text');\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nreturn \"done\";", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-10 Wednesday\", \"2024-04-11 Thursday\"),\n    @(\"482\u00f79=53, 5\", \"537\u00f73=179, 0\"),\n    @(\"738\u00f78=92, 2\", \"721\u00f77=103, 0\"),\n    @(\"741\u00f72=370, 1\", \"827\u00f75=165, 2\"),\n    @(\"931\u00f75=186, 1\", \"789\u00f72=394, 1\"),\n    @(\"893\u00f73=297, 2\", \"192\u00f73=64, 0\"),\n    @(\"877\u00f79=97, 4\", \"164\u00f77=23, 3\"),\n    @(\"501\u00f72=250, 1\", \"986\u00f78=123, 2\"),\n    @(\"317\u00f72=158, 1\", \"236\u00f79=26, 2\"),\n    @(\"351\u00f79=39, 0\", \"246\u00f76=41, 0\"),\n    @(\"419\u00f78=52, 3\", \"839\u00f77=119, 6\"),\n    @(\"882\u00f77=126, 0\", \"517\u00f74=129, 1\"),\n    @(\"846\u00f79=94, 0\", \"109\u00f75=21, 4\"),\n    @(\"338\u00f72=169, 0\", \"988\u00f79=109, 7\"),\n    @(\"895\u00f74=223, 3\", \"627\u00f76=104, 3\"),\n    @(\"781\u00f78=97, 5\", \"896\u00f79=99, 5\"),\n    @(\"257\u00f75=51, 2\", \"233\u00f74=58, 1\"),\n    @(\"282\u00f74=70, 2\", \"110\u00f78=13, 6\"),\n    @(\"645\u00f73=215, 0\", \"182\u00f79=20, 2\"),\n    @(\"496\u00f78=62, 0\", \"140\u00f76=23, 2\"),\n    @(\"785\u00f79=87, 2\", \"478\u00f75=95, 3\"),\n    @(\"636\u00f77=90, 6\", \"935\u00f77=133, 4\"),\n    @(\"771\u00f72=385, 1\", \"693\u00f73=231, 0\"),\n    @(\"877\u00f72=438, 1\", \"785\u00f79=87, 2\"),\n    @(\"912\u00f79=101, 3\", \"756\u00f78=94, 4\"),\n    @(\"914\u00f72=457, 0\", \"680\u00f74=170, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n}"}
